# Weekly manga ranking update: add the "2026-01-14" sheet.
#
# The simplest way to get an exact structural match (sheetPr/outline props,
# sheetFormatPr, pageMargins, header style, etc.) is to duplicate the most
# recent weekly sheet ("2026-01-07") -- which already carries the right
# sheet-level properties and the bold/bordered header style plus the
# light-yellow "new volume" highlight style -- rename the copy, and then
# overwrite its data rows with this week's ranking. Per-row highlight
# formatting is reset explicitly for every data row (cleared, or (re)applied
# by copying the format from a cell already using the highlight style) so no
# stale formatting from the previous week's sheet survives the overwrite.

$wb = $excel.ActiveWorkbook

$templateName = "2026-01-07"
$newName = "2026-01-14"

$templateSheet = $wb.Worksheets.Item($templateName)
$templateSheet.Copy([System.Reflection.Missing]::Value, $templateSheet)

# The duplicate is inserted right after the template and named "<template> (2)".
$ws = $wb.Worksheets.Item("$templateName (2)")
$ws.Name = $newName

# A cell already carrying the light-yellow "new volume" highlight style, used
# as a format-paint source so the duplicated sheet's existing style table
# entry is reused instead of a new one being minted. Pulled from the
# (untouched) template sheet so it stays valid for the whole loop below.
$highlightSource = $templateSheet.Range("C17")

# --- Ranking rows: rank, title, volume, highlighted? (1 = light-yellow "new volume" fill) ---
$rows = @(
    @(1, 'ダンダダン', 22, 0),
    @(2, 'カグラバチ', 10, 0),
    @(3, '呪術廻戦≡(モジュロ)', 1, 1),
    @(4, '片田舎のおっさん、剣聖になる~ただの田舎の剣術師範だったのに、大成した弟子たちが俺を放ってくれない件~', 8, 0),
    @(5, 'るろうに剣心―明治剣客浪漫譚・北海道編―', 10, 0),
    @(6, '葬送のフリーレン', 15, 0),
    @(7, '冒険王ビィト', 19, 0),
    @(8, '人間カード', 2, 1),
    @(9, '俺の死亡フラグが留まるところを知らない', 2, 1),
    @(10, 'ワスレモノ', 1, 1),
    @(11, '朱にまじわれば', 1, 1),
    @(12, '年の差十五の旦那様~辺境伯の花嫁候補~', 6, 0),
    @(13, '人間カード', 1, 1),
    @(14, '葬-はぶり-', 1, 1),
    @(15, '追放されたので、前世のレシピでパン職人はじめます!1', 1, 1),
    @(16, '追放された転生重騎士はゲーム知識で無双する', 16, 0),
    @(17, 'あかね噺', 20, 0),
    @(18, 'リアル炎上「GPS」', 1, 1),
    @(19, 'カムゴロシ', 3, 1),
    @(20, '隔離都市', 1, 1),
    @(21, '年の差十五の旦那様~辺境伯の花嫁候補~', 5, 0),
    @(22, '俺だけレベルアップな件', 23, 0),
    @(23, '俺んちに来た女騎士と田舎暮らしすることになった件', 2, 1),
    @(24, '年の差十五の旦那様~辺境伯の花嫁候補~', 4, 0),
    @(25, 'ゴブリンスレイヤー', 17, 0),
    @(26, '俺んちに来た女騎士と田舎暮らしすることになった件', 1, 1),
    @(27, '邪風のストラ', 1, 1),
    @(28, '僕たちの生きた理由', 1, 1),
    @(29, '忘却バッテリー', 23, 0),
    @(30, 'ドラゴンクエスト ダイの大冒険 勇者アバンと獄炎の魔王', 14, 0),
    @(31, 'キン肉マン', 91, 0),
    @(32, '転生先で推しの弟に美味しくいただかれています', 4, 0),
    @(33, 'カムゴロシ', 1, 1),
    @(34, 'カムゴロシ', 2, 1),
    @(35, '「才能の器」で目指す迷宮最深部 スキル横伸ばしのはずが、万能チートだった!', 2, 1),
    @(36, 'ガチャを回して仲間を増やす 最強の美少女軍団を作り上げろ THE COMIC', 5, 0),
    @(37, '新テニスの王子様', 46, 0),
    @(38, 'ハニーレモンソーダ', 30, 0),
    @(39, '脱獄のカザリヤ', 2, 1),
    @(40, 'ステータス・オール∞(インフィニティ) ∞使いの最強能力者、異世界を自由気ままに暮らします!', 1, 1),
    @(41, '次に買うマンガ、この1話で決めよう! ~異世界デビューはここから! 異世界初心者編~', 2, 1),
    @(42, '365日前の花嫁~男友達から強引に求婚されています', 5, 0),
    @(43, '結婚式当日、新郎の弟にプロポーズされました。', 6, 0),
    @(44, '魔法歌姫マジカルギンガ 第26話', 26, 0),
    @(45, '猫と竜', 2, 1),
    @(46, '「才能の器」で目指す迷宮最深部 スキル横伸ばしのはずが、万能チートだった!', 1, 1),
    @(47, '不遇職の成り上がり 美少女人形と最強まで最高速で上りつめる', 1, 1),
    @(48, 'ステータス・オール∞(インフィニティ) ∞使いの最強能力者、異世界を自由気ままに暮らします!', 2, 1),
    @(49, 'アイヲンモール異世界店、本日グランドオープン! THE COMIC', 1, 1),
    @(50, 'ガチャを回して仲間を増やす 最強の美少女軍団を作り上げろ THE COMIC', 3, 1),
    @(51, 'ガチャを回して仲間を増やす 最強の美少女軍団を作り上げろ THE COMIC', 4, 0),
    @(52, 'ゲーム中盤で死ぬ悪役貴族に転生したので、外れスキルを駆使して最強を目指してみた', 4, 0),
    @(53, 'ワンパンマン', 35, 0),
    @(54, 'ONE PIECE', 113, 0),
    @(55, '脱獄のカザリヤ', 1, 1),
    @(56, 'イジメの時間', 1, 1),
    @(57, '俺の死亡フラグが留まるところを知らない', 1, 1),
    @(58, '悪魔なボクは退魔師サマに愛されたい!!', 1, 1),
    @(59, '突然パパになった最強ドラゴンの子育て日記~かわいい娘、ほのぼのと人間界最強に育つ~ THE COMIC', 1, 1),
    @(60, '規格外のダンジョン攻略者、実は異世界帰りの元勇者1', 1, 1),
    @(61, '偽装カレシに愛されてしまいました', 3, 1),
    @(62, 'くれなゐの花嫁~大正北國恋物語~', 1, 1),
    @(63, '男友達が激甘カレシになりました', 3, 1),
    @(64, 'What Happened to Yukio Sukimoto', 1, 1),
    @(65, '小悪魔くんの甘い囁き', 8, 0),
    @(66, '愛になるまであたためて', 9, 0),
    @(67, '転生先で推しの弟に美味しくいただかれています', 3, 1),
    @(68, '雑用付与術師が自分の最強に気付くまで(コミック)', 10, 0),
    @(69, '薬屋のひとりごと', 16, 0),
    @(70, '薬屋のひとりごと~猫猫の後宮謎解き手帳~', 21, 0),
    @(71, 'ケントゥリア', 7, 0),
    @(72, 'チェンソーマン', 21, 0),
    @(73, '猫と竜', 1, 1),
    @(74, '竜使の花嫁 ~新緑の乙女は聖竜の守護者に愛される~ 1(アリアンローズコミックス)', 1, 1),
    @(75, 'イジメの時間', 2, 1),
    @(76, '賢者の弟子を名乗る賢者 THE COMIC', 2, 1),
    @(77, '賢者の弟子を名乗る賢者 THE COMIC', 3, 1),
    @(78, '賢者の弟子を名乗る賢者 THE COMIC', 4, 0),
    @(79, '経験人数が見えるメガネ', 1, 1),
    @(80, '修羅幼女の英雄譚~半端者と言われた傭兵、幼女に転生して成り上がる~1', 1, 1),
    @(81, '偽装カレシに愛されてしまいました', 1, 1),
    @(82, '偽装カレシに愛されてしまいました', 2, 1),
    @(83, '俺様婚約者には惚れたくありません!@COMIC 第1話', 1, 1),
    @(84, '大正浪漫 斜陽のくちづけ~傷だらけのご令嬢は剛腕社長に一途に愛される', 1, 1),
    @(85, '暴君皇子の契約妃', 2, 1),
    @(86, '暴君皇子の契約妃', 3, 1),
    @(87, '伯爵家の不幸な養女は、異国の王子に愛される', 4, 0),
    @(88, '男友達が激甘カレシになりました', 1, 1),
    @(89, '男友達が激甘カレシになりました', 2, 1),
    @(90, '31歳、初カレ。', 2, 1),
    @(91, '結婚式当日、新郎の弟にプロポーズされました。', 5, 0),
    @(92, '年の差十五の旦那様~辺境伯の花嫁候補~', 2, 1),
    @(93, '年の差十五の旦那様~辺境伯の花嫁候補~', 3, 1),
    @(94, '365日前の花嫁~男友達から強引に求婚されています', 4, 0),
    @(95, 'ガチャを回して仲間を増やす 最強の美少女軍団を作り上げろ THE COMIC', 2, 1),
    @(96, '神血の救世主~0.00000001%を引き当て最強へ~', 9, 0),
    @(97, '神血の救世主~0.00000001%を引き当て最強へ~', 10, 0),
    @(98, '幼稚園WARS', 16, 0),
    @(99, '彼岸島 48日後…', 52, 0),
    @(100, 'みいちゃんと山田さん', 5, 0)
)

foreach ($row in $rows) {
    $r = $row[0] + 1
    $rankCell = $ws.Cells.Item($r, 1)
    $titleCell = $ws.Cells.Item($r, 2)
    $volumeCell = $ws.Cells.Item($r, 3)

    $rankCell.Value = $row[0]
    $titleCell.Value = $row[1]
    $volumeCell.ClearFormats()
    $volumeCell.Value = $row[2]

    if ($row[3] -eq 1) {
        $highlightSource.Copy()
        $volumeCell.PasteSpecial(-4122)
    }
}

Write-Output "Added sheet $newName with $($rows.Count) ranking rows."
